# Add 2022-Q3 data
# ------------------------------------------------------------------
# 1) "总计" sheet: insert a new leading data row for 2022-Q3 and
#    shift the existing quarters down by one row.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# Push existing rows 2..6 down to 3..7 (copy B:D values+formats, bottom-up
# so we never overwrite a row before reading it). Column A keeps its
# original per-row value/style; only the brand-new row 7 needs a fresh A7.
$summary.Range("A6:D6").Copy($summary.Range("A7:D7"))
$summary.Range("B5:D5").Copy($summary.Range("B6:D6"))
$summary.Range("B4:D4").Copy($summary.Range("B5:D5"))
$summary.Range("B3:D3").Copy($summary.Range("B4:D4"))
$summary.Range("B2:D2").Copy($summary.Range("B3:D3"))

$summary.Range("A7").Value = 5

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 8
$summary.Range("D2").Value = 0.32

# ------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" worksheet right after "总计" (i.e.
#    before the current second sheet, "2022-Q2") holding the fund
#    holdings detail for the quarter.
# ------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Header row - copy formatting (bold/border/center, style index 2) from
# an existing header cell on "总计" so the new sheet matches house style.
$summary.Range("B1").Copy($q3.Range("B1"))
$q3.Range("B1").Value = "基金代码"
$summary.Range("B1").Copy($q3.Range("C1"))
$q3.Range("C1").Value = "基金名称"
$summary.Range("B1").Copy($q3.Range("D1"))
$q3.Range("D1").Value = "基金规模"
$summary.Range("B1").Copy($q3.Range("E1"))
$q3.Range("E1").Value = "股票总仓位"
$summary.Range("B1").Copy($q3.Range("F1"))
$q3.Range("F1").Value = "仓位占比"
$summary.Range("B1").Copy($q3.Range("G1"))
$q3.Range("G1").Value = "持有市值(亿元)"
$summary.Range("B1").Copy($q3.Range("H1"))
$q3.Range("H1").Value = "仓位排名"

# Data rows. Column A (row index) reuses the "总计" A-column style
# (style index 2: bold/border/center). B/C columns are plain text.
# D/E/F columns are text-formatted numbers (kept as strings, matching
# the source data export). G is usually a text-formatted number too,
# except where the underlying market value is a literal zero - there
# it is stored as a true numeric 0 (matches the official export quirk
# also seen on the "2022-Q2" sheet's last row). H is a plain integer.

function Set-Row($r, $code, $name, $scale, $pos, $ratio, $mv, $mvIsNumber, $rank) {
    $summary.Range("A2").Copy($q3.Range("A$r"))
    $q3.Range("A$r").Value = $r - 2

    $q3.Range("B$r").NumberFormat = "@"
    $q3.Range("B$r").Value = $code
    $q3.Range("C$r").Value = $name

    $q3.Range("D$r").NumberFormat = "@"
    $q3.Range("D$r").Value = $scale

    $q3.Range("E$r").NumberFormat = "@"
    $q3.Range("E$r").Value = $pos

    $q3.Range("F$r").NumberFormat = "@"
    $q3.Range("F$r").Value = $ratio

    if ($mvIsNumber) {
        $q3.Range("G$r").Value = 0
    } else {
        $q3.Range("G$r").NumberFormat = "@"
        $q3.Range("G$r").Value = $mv
    }

    $q3.Range("H$r").Value = $rank
}

Set-Row 2 "006648" "汇安多因子混合A"         "3.69" "93.49" "2.70" "0.0996" $false 7
Set-Row 3 "001403" "招商国企改革主题混合"     "1.97" "85.29" "4.23" "0.0833" $false 10
Set-Row 4 "006649" "汇安多因子混合C"         "1.95" "93.49" "2.70" "0.0526" $false 7
Set-Row 5 "010558" "汇安鑫利优选混合A"       "1.30" "93.99" "2.62" "0.0341" $false 9
Set-Row 6 "002271" "招商安弘灵活配置混合"     "0.55" "53.72" "4.15" "0.0228" $false 7
Set-Row 7 "010559" "汇安鑫利优选混合C"       "0.65" "93.99" "2.62" "0.0170" $false 9
Set-Row 8 "011054" "申万菱信安鑫智选混合A"   "0.83" "22.80" "1.03" "0.0085" $false 7
Set-Row 9 "011055" "申万菱信安鑫智选混合C"   "0.00" "22.80" "1.03" "0"      $true  7

# Restore the originally-active sheet ("2021-Q2", the last tab) so the
# workbook's active-tab/selection state is left as untouched as possible,
# matching the diff (which does not touch bookViews/sheetView selection).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

Write-Output "2022-Q3 sheet populated"
